# Update the "dSF" (column F) values for specific rows to reflect the
# repulled / recalculated data, per commit: "repull data, push all data,
# mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F
$updates = @{
    4  = -3
    6  = -1
    11 = -2
    13 = -4
    16 = 0
    21 = 0
    24 = -2
    25 = -3
    28 = 0
    32 = -1
    34 = 1
    39 = -1
    41 = 4
    42 = -7
    50 = 5
    58 = -2
    62 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
